$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 231
$ws.Range("I12").Value = 240
$ws.Range("J12").Value = 217.5
$ws.Range("K12").Value = 240
$ws.Range("L12").Value = 217.5
$ws.Range("M12").Value = -70
$ws.Range("N12").Value = -557.5
$ws.Range("H39").Value = 57.22222
$ws.Range("I39").Value = 41.42857
$ws.Range("K39").Value = 124.28571
$ws.Range("M39").Value = 171.71429
$ws.Range("H43").Value = 13143.895
$ws.Range("I43").Value = 2281.1538
$ws.Range("K43").Value = 2281.1538
$ws.Range("M43").Value = -2212.1538
$ws.Range("H62").Value = 7499.75
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 7499
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 7499
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -8747
$ws.Range("H65").Value = 7499.75
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 7499
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 37495
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -43735
$ws.Range("H80").Value = 875.25
$ws.Range("I80").Value = 733
$ws.Range("K80").Value = 2199
$ws.Range("M80").Value = -1201
$ws.Range("H83").Value = 875.25
$ws.Range("I83").Value = 733
$ws.Range("K83").Value = 6597
$ws.Range("M83").Value = -1605
$ws.Range("H86").Value = 2023.5
$ws.Range("I86").Value = 1550
$ws.Range("J86").Value = 2497
$ws.Range("K86").Value = 1550
$ws.Range("L86").Value = 2497
$ws.Range("M86").Value = -427
$ws.Range("N86").Value = -4743
$ws.Range("H89").Value = 2023.5
$ws.Range("I89").Value = 1550
$ws.Range("J89").Value = 2497
$ws.Range("K89").Value = 7750
$ws.Range("L89").Value = 12485
$ws.Range("M89").Value = -2134
$ws.Range("N89").Value = -23717
$ws.Range("H92").Value = 1739.1333
$ws.Range("I92").Value = 507.66666
$ws.Range("J92").Value = 6665
$ws.Range("K92").Value = 507.66666
$ws.Range("L92").Value = 6665
$ws.Range("M92").Value = 740.33334
$ws.Range("N92").Value = -9161
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H107").Value = 492.83334
$ws.Range("I107").Value = 374.35
$ws.Range("K107").Value = 374.35
$ws.Range("M107").Value = 1545.65
$ws.Range("H108").Value = 56262.5
$ws.Range("I108").Value = 50525
$ws.Range("K108").Value = 50525
$ws.Range("M108").Value = -46685
$ws.Range("H132").Value = 1669272.4
$ws.Range("I132").Value = 3196.6667
$ws.Range("K132").Value = 9590.000100000001
$ws.Range("M132").Value = -7060.000100000001
$ws.Range("H137").Value = 2422.818
$ws.Range("I137").Value = 1912.1666
$ws.Range("K137").Value = 5736.4998
$ws.Range("M137").Value = -3186.4998
$ws.Range("H138").Value = 3322.275
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3322.275
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 9966.825000000001
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20246.825
$ws.Range("H141").Value = 5109.231
$ws.Range("I141").Value = 4115.75
$ws.Range("J141").Value = 6698.8
$ws.Range("K141").Value = 12347.25
$ws.Range("L141").Value = 20096.4
$ws.Range("M141").Value = -7167.25
$ws.Range("N141").Value = -30456.4

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1250.1428
$ws.Range("I2").Value = 926.7059
$ws.Range("K2").Value = 926.7059
$ws.Range("M2").Value = -813.7059
$ws.Range("H32").Value = 17361.555
$ws.Range("I32").Value = 6575.0464
$ws.Range("J32").Value = 59527
$ws.Range("K32").Value = 6575.0464
$ws.Range("L32").Value = 59527
$ws.Range("M32").Value = -6288.0464
$ws.Range("N32").Value = -60101
$ws.Range("H61").Value = 4973.273
$ws.Range("I61").Value = 4952.5
$ws.Range("J61").Value = 4985.143
$ws.Range("K61").Value = 4952.5
$ws.Range("L61").Value = 4985.143
$ws.Range("M61").Value = -4740.5
$ws.Range("N61").Value = -5409.143
$ws.Range("H96").Value = 14159.6
$ws.Range("J96").Value = 14159.6
$ws.Range("L96").Value = 14159.6
$ws.Range("N96").Value = -19651.6
$ws.Range("H97").Value = 1109.0286
$ws.Range("I97").Value = 1070.8276
$ws.Range("J97").Value = 1293.6666
$ws.Range("K97").Value = 1070.8276
$ws.Range("L97").Value = 1293.6666
$ws.Range("M97").Value = -574.8276000000001
$ws.Range("N97").Value = -2285.6666
$ws.Range("H102").Value = 2508.4075
$ws.Range("I102").Value = 2515.0417
$ws.Range("J102").Value = 2455.3333
$ws.Range("K102").Value = 2515.0417
$ws.Range("L102").Value = 2455.3333
$ws.Range("M102").Value = -893.0417000000002
$ws.Range("N102").Value = -5699.3333
$ws.Range("H110").Value = 2792.75
$ws.Range("I110").Value = 2391.3125
$ws.Range("J110").Value = 4398.5
$ws.Range("K110").Value = 2391.3125
$ws.Range("L110").Value = 4398.5
$ws.Range("M110").Value = -346.3125
$ws.Range("N110").Value = -8488.5
$ws.Range("H116").Value = 1250.1428
$ws.Range("I116").Value = 926.7059
$ws.Range("K116").Value = 926.7059
$ws.Range("M116").Value = 1367.2941
$ws.Range("H132").Value = 1499.2858
$ws.Range("I132").Value = 849.4167
$ws.Range("K132").Value = 2548.2501
$ws.Range("M132").Value = -18.2501000000002
$ws.Range("H136").Value = 4973.273
$ws.Range("I136").Value = 4952.5
$ws.Range("J136").Value = 4985.143
$ws.Range("K136").Value = 14857.5
$ws.Range("L136").Value = 14955.429
$ws.Range("M136").Value = -12307.5
$ws.Range("N136").Value = -20055.429

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1250.1428
$ws.Range("I3").Value = 926.7059
$ws.Range("K3").Value = 926.7059
$ws.Range("M3").Value = -812.7059
$ws.Range("H22").Value = 659.4091
$ws.Range("I22").Value = 553.41174
$ws.Range("K22").Value = 553.41174
$ws.Range("M22").Value = -380.41174
$ws.Range("H80").Value = 1138.5
$ws.Range("I80").Value = 250
$ws.Range("J80").Value = 2027
$ws.Range("K80").Value = 250
$ws.Range("L80").Value = 2027
$ws.Range("M80").Value = 748
$ws.Range("N80").Value = -4023
$ws.Range("H83").Value = 1138.5
$ws.Range("I83").Value = 250
$ws.Range("J83").Value = 2027
$ws.Range("K83").Value = 1250
$ws.Range("L83").Value = 10135
$ws.Range("M83").Value = 3742
$ws.Range("N83").Value = -20119
$ws.Range("H99").Value = 1518.5385
$ws.Range("I99").Value = 1494.2
$ws.Range("K99").Value = 1494.2
$ws.Range("M99").Value = 3.799999999999955
$ws.Range("H105").Value = 3248.75
$ws.Range("I105").Value = 2828.9
$ws.Range("J105").Value = 3948.5
$ws.Range("K105").Value = 2828.9
$ws.Range("L105").Value = 3948.5
$ws.Range("M105").Value = -1081.9
$ws.Range("N105").Value = -7442.5
$ws.Range("H107").Value = 1368.25
$ws.Range("I107").Value = 1109.1666
$ws.Range("K107").Value = 1109.1666
$ws.Range("M107").Value = 810.8334
$ws.Range("H134").Value = 2886.6538
$ws.Range("I134").Value = 2886.6538
$ws.Range("K134").Value = 8659.9614
$ws.Range("M134").Value = -6124.9614

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H3").Value = 11000
$ws.Range("J3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15226
$ws.Range("H16").Value = 1822
$ws.Range("I16").Value = 2349
$ws.Range("K16").Value = 2349
$ws.Range("M16").Value = -2062
$ws.Range("H22").Value = 195
$ws.Range("I22").Value = 195
$ws.Range("K22").Value = 195
$ws.Range("M22").Value = 155
$ws.Range("H31").Value = 6490.5
$ws.Range("I31").Value = 3595.9143
$ws.Range("J31").Value = 10895.305
$ws.Range("K31").Value = 3595.9143
$ws.Range("L31").Value = 10895.305
$ws.Range("M31").Value = -3300.9143
$ws.Range("N31").Value = -11485.305
$ws.Range("H34").Value = 6490.5
$ws.Range("I34").Value = 3595.9143
$ws.Range("J34").Value = 10895.305
$ws.Range("K34").Value = 3595.9143
$ws.Range("L34").Value = 10895.305
$ws.Range("M34").Value = -3393.9143
$ws.Range("N34").Value = -11299.305
$ws.Range("H113").Value = 1822
$ws.Range("I113").Value = 2349
$ws.Range("K113").Value = 2349
$ws.Range("M113").Value = -179
$ws.Range("H132").Value = 2371.889
$ws.Range("I132").Value = 2177.64
$ws.Range("K132").Value = 6532.92
$ws.Range("M132").Value = -4002.92
$ws.Range("H134").Value = 975.3913
$ws.Range("I134").Value = 953.3333
$ws.Range("K134").Value = 2859.9999
$ws.Range("M134").Value = -324.9998999999998

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 41.166668
$ws.Range("I2").Value = 41.166668
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 247.000008
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -134.000008
$ws.Range("N2").ClearContents()
$ws.Range("H23").Value = 127.22222
$ws.Range("J23").Value = 129
$ws.Range("L23").Value = 387
$ws.Range("N23").Value = -857
$ws.Range("H29").Value = 558.8333
$ws.Range("I29").Value = 749.5
$ws.Range("J29").Value = 177.5
$ws.Range("K29").Value = 2248.5
$ws.Range("L29").Value = 532.5
$ws.Range("M29").Value = -1971.5
$ws.Range("N29").Value = -1086.5
$ws.Range("H55").Value = 1883.5714
$ws.Range("J55").Value = 3666.6667
$ws.Range("L55").Value = 11000.0001
$ws.Range("N55").Value = -11354.0001
$ws.Range("H92").Value = 309.5
$ws.Range("I92").Value = 370
$ws.Range("K92").Value = 1110
$ws.Range("M92").Value = 138
$ws.Range("H129").Value = 7518.769
$ws.Range("J129").Value = 9179.5
$ws.Range("L129").Value = 27538.5
$ws.Range("N129").Value = -37538.5
$ws.Range("H131").Value = 55359.953
$ws.Range("J131").Value = 2680.5557
$ws.Range("L131").Value = 8041.6671
$ws.Range("N131").Value = -18121.6671
$ws.Range("H132").Value = 1281
$ws.Range("I132").Value = 1338.875
$ws.Range("J132").Value = 1049.5
$ws.Range("K132").Value = 12049.875
$ws.Range("L132").Value = 9445.5
$ws.Range("M132").Value = -9519.875
$ws.Range("N132").Value = -14505.5
$ws.Range("H137").Value = 18892.666
$ws.Range("I137").Value = 18892.666
$ws.Range("K137").Value = 56677.99800000001
$ws.Range("M137").Value = -51577.99800000001
$ws.Range("H139").Value = 10586.5
$ws.Range("I139").Value = 8403.571
$ws.Range("J139").Value = 13133.25
$ws.Range("K139").Value = 25210.713
$ws.Range("L139").Value = 39399.75
$ws.Range("M139").Value = -20070.713
$ws.Range("N139").Value = -49679.75

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 121.86957
$ws.Range("I2").Value = 136.52632
$ws.Range("J2").Value = 52.25
$ws.Range("K2").Value = 136.52632
$ws.Range("L2").Value = 52.25
$ws.Range("M2").Value = -23.52632
$ws.Range("N2").Value = -278.25
$ws.Range("H9").Value = 289.25
$ws.Range("I9").Value = 52.333332
$ws.Range("K9").Value = 52.333332
$ws.Range("M9").Value = 117.666668
$ws.Range("H24").Value = 32492
$ws.Range("I24").Value = 30000
$ws.Range("J24").Value = 33738
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 33738
$ws.Range("M24").Value = -29827
$ws.Range("N24").Value = -34084
$ws.Range("H46").Value = 5225
$ws.Range("I46").Value = 5225
$ws.Range("K46").Value = 5225
$ws.Range("M46").Value = -5069
$ws.Range("H70").Value = 7688.25
$ws.Range("I70").Value = 5008
$ws.Range("K70").Value = 5008
$ws.Range("M70").Value = -4738
$ws.Range("H73").Value = 7688.25
$ws.Range("I73").Value = 5008
$ws.Range("K73").Value = 5008
$ws.Range("M73").Value = -4072
$ws.Range("H80").Value = 10281.909
$ws.Range("I80").Value = 2800.5715
$ws.Range("K80").Value = 2800.5715
$ws.Range("M80").Value = -1802.5715
$ws.Range("H83").Value = 10281.909
$ws.Range("I83").Value = 2800.5715
$ws.Range("K83").Value = 14002.8575
$ws.Range("M83").Value = -9010.8575
$ws.Range("H97").Value = 1545.2941
$ws.Range("I97").Value = 1520.4286
$ws.Range("K97").Value = 1520.4286
$ws.Range("M97").Value = -1024.4286
$ws.Range("H107").Value = 869.8
$ws.Range("I107").Value = 568.0909
$ws.Range("J107").Value = 1699.5
$ws.Range("K107").Value = 568.0909
$ws.Range("L107").Value = 1699.5
$ws.Range("M107").Value = 1351.9091
$ws.Range("N107").Value = -5539.5
$ws.Range("H122").Value = 57270.277
$ws.Range("I122").Value = 201141.2
$ws.Range("J122").Value = 1935.3077
$ws.Range("K122").Value = 603423.6000000001
$ws.Range("L122").Value = 5805.9231
$ws.Range("M122").Value = -600973.6000000001
$ws.Range("N122").Value = -10705.9231
$ws.Range("H132").Value = 2173.4688
$ws.Range("I132").Value = 2153.5
$ws.Range("K132").Value = 6460.5
$ws.Range("M132").Value = -3930.5

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H50").Value = 34799.6
$ws.Range("I50").Value = 37999
$ws.Range("J50").Value = 32666.666
$ws.Range("K50").Value = 37999
$ws.Range("L50").Value = 32666.666
$ws.Range("M50").Value = -37362
$ws.Range("N50").Value = -33940.666
$ws.Range("H61").Value = 1315.8
$ws.Range("J61").Value = 1050
$ws.Range("L61").Value = 1050
$ws.Range("N61").Value = -1454
$ws.Range("H68").Value = 4150
$ws.Range("I68").Value = 5500
$ws.Range("K68").Value = 5500
$ws.Range("M68").Value = -4751
$ws.Range("H71").Value = 4150
$ws.Range("I71").Value = 5500
$ws.Range("K71").Value = 27500
$ws.Range("M71").Value = -23756
$ws.Range("H100").Value = 3499.1428
$ws.Range("I100").Value = 4099.75
$ws.Range("K100").Value = 4099.75
$ws.Range("M100").Value = -3558.75
$ws.Range("H113").Value = 1315.8
$ws.Range("J113").Value = 1050
$ws.Range("L113").Value = 1050
$ws.Range("N113").Value = -5390
$ws.Range("H116").Value = 61998.8
$ws.Range("J116").Value = 59998.75
$ws.Range("L116").Value = 59998.75
$ws.Range("N116").Value = -69176.75
$ws.Range("H122").Value = 3198.842
$ws.Range("I122").Value = 1746.5
$ws.Range("J122").Value = 3869.1538
$ws.Range("K122").Value = 5239.5
$ws.Range("L122").Value = 11607.4614
$ws.Range("M122").Value = -2789.5
$ws.Range("N122").Value = -16507.4614
$ws.Range("H136").Value = 3541.258
$ws.Range("I136").Value = 2949
$ws.Range("K136").Value = 8847
$ws.Range("M136").Value = -6297

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H11").Value = 31599
$ws.Range("J11").Value = 32398.5
$ws.Range("L11").Value = 32398.5
$ws.Range("N11").Value = -32682.5
$ws.Range("H51").Value = 16233.333
$ws.Range("J51").Value = 17442.857
$ws.Range("L51").Value = 17442.857
$ws.Range("N51").Value = -18462.857
$ws.Range("H62").Value = 18487.5
$ws.Range("J62").Value = 18487.5
$ws.Range("L62").Value = 18487.5
$ws.Range("N62").Value = -19735.5
$ws.Range("H65").Value = 18487.5
$ws.Range("J65").Value = 18487.5
$ws.Range("L65").Value = 92437.5
$ws.Range("N65").Value = -98677.5
$ws.Range("H81").Value = 2950
$ws.Range("I81").Value = 900
$ws.Range("K81").Value = 1800
$ws.Range("M81").Value = -739
$ws.Range("H84").Value = 2950
$ws.Range("I84").Value = 900
$ws.Range("K84").Value = 9000
$ws.Range("M84").Value = -3696
$ws.Range("H107").Value = 851.9524
$ws.Range("I107").Value = 609.61536
$ws.Range("J107").Value = 1245.75
$ws.Range("K107").Value = 1828.84608
$ws.Range("L107").Value = 3737.25
$ws.Range("M107").Value = 91.15391999999997
$ws.Range("N107").Value = -7577.25
$ws.Range("H113").Value = 429.4
$ws.Range("I113").Value = 429.4
$ws.Range("K113").Value = 1288.2
$ws.Range("M113").Value = 881.8000000000002
$ws.Range("H122").Value = 10031.5
$ws.Range("I122").Value = 9000.689
$ws.Range("J122").Value = 19996
$ws.Range("K122").Value = 27002.067
$ws.Range("L122").Value = 59988
$ws.Range("M122").Value = -24552.067
$ws.Range("N122").Value = -64888
$ws.Range("H136").Value = 2954.0833
$ws.Range("I136").Value = 3004.9
$ws.Range("K136").Value = 9014.700000000001
$ws.Range("M136").Value = -6464.700000000001
